$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Majorelle Magdy, Dr. Eman Tantawi'
$ws.Range("G3").Value = 'Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Asmaa Reda'
$ws.Range("G4").Value = 'Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud'
$ws.Range("G5").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Eman Tantawi, Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Hend Mahmoud'
$ws.Range("G6").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Hend Mahmoud'
$ws.Range("G7").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Hend Mahmoud'
$ws.Range("G8").Value = 'Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Administrator, Dr. Manar Montaser, Dr. Asmaa Reda'
$ws.Range("G9").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Gehan Adel, Dr. Manar Montaser, Dr. Asmaa Reda, Dr. Hend Mahmoud'
$ws.Range("G10").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Heba Mahmoud Ali, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Gehan Adel, Dr. Sara Wael, Dr. Alshimaa Atef'
$ws.Range("G11").Value = 'Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Asmaa Reda'
$ws.Range("G13").Value = 'Dr. Omnia Mohammad, Dr. Shimaa Ashraf, Dr. Safa Hany, Dr. Mariam Nour El-Din, D Wessam Atef'
$ws.Range("G14").Value = 'Dr. Shimaa Ashraf, Dr. Safa Hany'
$ws.Range("G15").Value = 'D Wessam Atef, Dr. Amal Awwad'
$ws.Range("G17").Value = 'Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Marwa Mustafa, Dr. Arwa Al-Sayed, Dr. Basma Hamed, Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Sarah Abdelmohsen'
$ws.Range("G19").Value = 'Dr. Sarah Mahdy, D Mariam E. Mohammad'
$ws.Range("G22").Value = 'Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed'
$ws.Range("G23").Value = 'Dr. Hana Amr, Dr. Nourham Mostafa'
$ws.Range("G24").Value = 'Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Marina Atef, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Remon'
$ws.Range("G25").Value = 'Dr. Marina Atef, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Remon'
$ws.Range("G27").Value = 'Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Eman Mohammad Al, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Remon'
$ws.Range("G28").Value = 'Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Abdullah El-Agrody, Dr. Remon'
$ws.Range("G29").Value = 'Dr. Neveen Nashaat, Dr. Naema Gomaa, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Remon'
$ws.Range("G30").Value = 'Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud'
$ws.Range("G31").Value = 'Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Asmaa Reda'
$ws.Range("G32").Value = 'Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud'
$ws.Range("G33").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Mohammad El-Tanany, Dr. Nesma, Dr. Eman Tantawi, Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Hend Mahmoud'
$ws.Range("G34").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Nahla Nagiub, Dr. Nourhan Mahmoud, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Hend Mahmoud'
$ws.Range("G35").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Asmaa Reda, Dr. Hend Mahmoud'
$ws.Range("G36").Value = 'Dr. Majorelle Magdy, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Administrator, Dr. Manar Montaser, Dr. Asmaa Reda'
$ws.Range("G37").Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Majorelle Magdy, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Gehan Adel, Dr. Manar Montaser, Dr. Asmaa Reda, Dr. Hend Mahmoud'
$ws.Range("G38").Value = 'Dr. Servinaz Sayed Mohammad, Dr. Heba Mahmoud Ali, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Gehan Adel, Dr. Sara Wael, Dr. Alshimaa Atef'
$ws.Range("G39").Value = 'Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Asmaa Reda'
$ws.Range("G41").Value = 'Dr. Omnia Mohammad, Dr. Shimaa Ashraf, Dr. Safa Hany, Dr. Mariam Nour El-Din, D Wessam Atef'
$ws.Range("G42").Value = 'Dr. Shimaa Ashraf, Dr. Safa Hany'
$ws.Range("G43").Value = 'D Wessam Atef, Dr. Amal Awwad'
$ws.Range("G45").Value = 'Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Esraa Mostafa, Dr. Nourhan Osama, Dr. Marwa Mustafa, Dr. Arwa Al-Sayed, Dr. Basma Hamed, Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Sarah Abdelmohsen'
$ws.Range("G47").Value = 'Dr. Sarah Mahdy, D Mariam E. Mohammad'
$ws.Range("G50").Value = 'Dr. Nancy Abd Al-Shafy, Dr. Amr Saeed'
$ws.Range("G51").Value = 'Dr. Hana Amr, Dr. Nourham Mostafa'
$ws.Range("G52").Value = 'Dr. Maryam Ashraf, Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Marina Atef, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Remon'
$ws.Range("G53").Value = 'Dr. Marina Atef, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Abdullah El-Agrody, Dr. Youstina Magdy, Dr. Aya Emad, Dr. Remon'
$ws.Range("G55").Value = 'Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Yasmin, Dr. Wafaa Ebida, Dr. Eman Mohammad Al, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Remon'
$ws.Range("G56").Value = 'Dr. Neveen Nashaat, Dr. Salma Hassan, Dr. Aya Hanafy, Dr. Wafaa Ebida, Dr. Eman Samir Gabry, Dr. Nardine, Dr. Abdullah El-Agrody, Dr. Remon'
$ws.Range("G57").Value = 'Dr. Neveen Nashaat, Dr. Naema Gomaa, Dr. Eman Samir Gabry, Dr. Ola Abd Al-Fattah, Dr. Monica, Dr. Remon'
